$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map each chromosome label (column A) to its RefSeq accession (NCBI GRCh37/hg19 build),
# used to build the new "refseq" and "crispick" columns.
$chrAccession = @{
    "11" = "NC_000011.9"
    "8"  = "NC_000008.10"
    "16" = "NC_000016.9"
    "6"  = "NC_000006.11"
    "17" = "NC_000017.10"
    "19" = "NC_000019.9"
    "7"  = "NC_000007.13"
    "12" = "NC_000012.11"
    "10" = "NC_000010.10"
    "9"  = "NC_000009.11"
    "20" = "NC_000020.10"
    "14" = "NC_000014.8"
    "2"  = "NC_000002.11"
    "5"  = "NC_000005.9"
    "15" = "NC_000015.9"
    "21" = "NC_000021.8"
    "4"  = "NC_000004.11"
}

# Header row: append the two new columns after "SNP" (F)
$ws.Cells.Item(1, 7).Value = "refseq"
$ws.Cells.Item(1, 8).Value = "crispick"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $chrom = [string]$ws.Cells.Item($r, 1).Value2
    $start = $ws.Cells.Item($r, 2).Value2
    $end   = $ws.Cells.Item($r, 3).Value2

    $accession = $chrAccession[$chrom]

    $ws.Cells.Item($r, 7).Value = $accession
    $ws.Cells.Item($r, 8).Value = "$($accession):+:$($start)-$($end)"
}
